# Generate Report for Handoff
# Updates the "Latest Handoff Datetime" timestamps and "Priority" values
# for the rows that were (re)processed during this handoff run, across
# the Overview / zh-cn / de-de worksheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Rows affected by this handoff generation (table data rows 8-11 and 13-14;
# row 12 was not part of this run).
$rows = @(8, 9, 10, 11, 13, 14)

foreach ($r in $rows) {
    # Overview sheet: column G = "Latest HO Xliff Generate Date"
    $overview.Range("G$r").Value = "2017-01-03 08:36:06"

    # zh-cn sheet: column E = "Priority", column H = "Latest Handoff Datetime"
    $zhcn.Range("E$r").Value = "ht"
    $zhcn.Range("H$r").Value = "2017-01-03 08:35:53"

    # de-de sheet: column E = "Priority", column H = "Latest Handoff Datetime"
    $dede.Range("E$r").Value = "ht"
    $dede.Range("H$r").Value = "2017-01-03 08:36:06"
}
